# Applies the edit described by the diff:
#  - A2: "Tablet" -> "HP ElitePad 1000 G2 Tablet"
#  - B5: (empty) -> "Funcionou"
#  - Selection moves to C7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "HP ElitePad 1000 G2 Tablet"
$ws.Range("B5").Value = "Funcionou"

$ws.Range("C7").Select()
